$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S (year 2021) values, one per data row, mirroring column R's
# formatting for that row (header / bold-data / regular-data / bottom-border row).
$values = [ordered]@{
    4  = 2021
    5  = 6.1
    6  = 1.6
    7  = 3.6
    8  = 27.2
    9  = 7.2
    10 = 2.6
    11 = 12.5
    12 = 6.4
    13 = 5.2
    14 = 0.9
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("S$row")
    $cell.Value = $values[$row]

    # Shared look across the whole block: right aligned, vertically centered,
    # Times New Roman 9pt (matches the rest of the data table).
    $cell.HorizontalAlignment = -4152   # xlRight
    $cell.VerticalAlignment = -4108     # xlCenter
    $cell.Font.Name = "Times New Roman"
    $cell.Font.Size = 9

    if ($row -eq 4) {
        # Header row: bold, general number format, boxed top+bottom medium border.
        $cell.NumberFormat = "General"
        $cell.Font.Bold = $true
        $cell.Borders.Item(8).LineStyle = 1
        $cell.Borders.Item(8).Weight = -4138
        $cell.Borders.Item(9).LineStyle = 1
        $cell.Borders.Item(9).Weight = -4138
    }
    elseif ($row -eq 5) {
        # First data row: bold, one decimal place, no border.
        $cell.NumberFormat = "0.0"
        $cell.Font.Bold = $true
    }
    elseif ($row -eq 14) {
        # Last data row: not bold, one decimal place, medium bottom border.
        $cell.NumberFormat = "0.0"
        $cell.Font.Bold = $false
        $cell.Borders.Item(9).LineStyle = 1
        $cell.Borders.Item(9).Weight = -4138
    }
    else {
        # Regular data rows: not bold, one decimal place, no border.
        $cell.NumberFormat = "0.0"
        $cell.Font.Bold = $false
    }
}

# Reflect the post-edit selection recorded for the sheet.
$ws.Range("Q19").Select()
